$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix Objetivos text in B10/C10 (row10 currently holds professor name by mistake)
$ws.Range("B10").Value2 = "Fornecer os conceitos básicos de Mecânica dos Fluidos e Transferência de Calor e Massa com aplicações à Engenharia. Possibilitar aos alunos uma base científica para que possam se desenvolver em demais disciplinas tecnológicas do curso."

# 2. Insert a new row at position 14 (shifts old rows 14-24 down to 15-25)
$ws.Rows.Item(14).Insert()

# 3. Row 13: drop the "Programa resumido:" label (A13) and put the professor name in B13/C13
$ws.Range("A13").Clear()
$ws.Range("B13").Value2 = "4808662 - Lucrécio Fábio dos Santos"

# 4. New row 14: "Programa resumido:" label + short syllabus text
$ws.Range("A14").Value2 = "Programa resumido:"
$ws.Range("B14").Value2 = "Conceitos ligados ao escoamento de fluídos e equações fundamentais, Escoamento incompressível de fluidos não viscosos, Escoamento viscoso incompressível, Transferência de Calor. Transferência de Massa"
$ws.Rows.Item(14).RowHeight = 60

# 5. Row 13 height reverts to default (no custom height)
$ws.Rows.Item(13).AutoFit()

# 6. Row 16 (was old row 15): "Programa:" gets the full syllabus text (was bogus date)
$ws.Range("B16").Value2 = "1. Conceitos ligados ao escoamento de fluídos e equações fundamentais1.1. Características e definições dos escoamentos;1.2. Conceitos de sistema e volume de controle;1.3. Equação da conservação da massa;1.4. Equação da conservação da energia;1.5. Equação da conservação da quantidade de movimento;1.6. Introdução à análise diferencial do movimento de fluidos.2. Escoamento incompressível de fluidos não viscosos2.1. Equação de Euler;2.2. Equação de Bernoulli;2.3. Aplicações da equação de Bernoulli.3. Escoamento viscoso incompressível3.1. Atrito e perda de carga;3.2. Avaliação das perdas de carga: regime laminar e turbulento;3.3. Equações de Hagen - Poiseuille e Darcy – Weisbach3.4. Diagrama de Moody e Moody –Rouse;3.5. Método dos comprimentos equivalentes.3.6. Presença de máquina no escoamento (bomba e turbina), Potência e rendimento;3.7. Medidores de vazão.4. Transferência de Calor4.1. Definição de Calor.4.2. Mecanismo da Condução.4.3. Mecanismo da Convecção.4.4. Associação de Mecanismos.5. Transferência de Massa5.1. Difusão e convecção mássica;5.2. 1ª lei de Fick;5.3. Concentrações mássica e molar;5.4. Frações mássica e molar;5.5. Velocidades médias mássica e molar;5.6. Fluxos difusivo mássico, difusivo molar, convectivo mássico e convectivo molar;5.7. Fluxo mássico total e fluxo molar total."

# 7. Row 19 (was old row 18): "Metodo:" gets the teaching-method text (was bogus professor name)
$ws.Range("B19").Value2 = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios, aulas de laboratório."

# 8. Row 20 (was old row 19): "Criterio:" gets the grading-formula text
$ws.Range("B20").Value2 = "Nota de duas provas (P1 e P2)Fórmula: M1 = (P1 + 2 x P2)/3.."

# 9. Row 21 (was old row 20): "Norma de recuperacao:" gets the recovery-exam text
$ws.Range("B21").Value2 = "Aplicação de uma prova envolvendo o assunto de todo semestre.NR (nota da recuperação) = (M1 + NR)/2."

# 10. Row 22 (was old row 21): "Bibliografia:" gets the bibliography text
$ws.Range("B22").Value2 = "1. FOX, R.W., MCDONALD, A.T., “Introdução à Mecânica dos Fluidos”, Ed. Guanabara Koogan.2. STREETER, V.L., WYLE,E.B., “Mecânica dos Fluidos”, Ed. Mc Graw Hill.3. OZISIK,M.N., “Transferência de Calor.”, Ed. Guanabara Koogan.4. INCROPERA, F.P.W., “Fundamentos de Transferência de Calor e Massa”, Ed. Guanabara Koogan.5. MUNSON, B.R.; YOUNG, D.F.; OKIISHI, T.H. Fundamentos da Mecânica dos Fluidos. Editora Edgard Blucher6 - GIORGETI, M. (2012) Fundamentos de Fenômenos de Transporte. Editora Campus"

# Mirror B-column values into C column for all the rows where the diff shows both B and C set
$ws.Range("C10").Value2 = $ws.Range("B10").Value2
$ws.Range("C13").Value2 = $ws.Range("B13").Value2
$ws.Range("C14").Value2 = $ws.Range("B14").Value2
$ws.Range("C16").Value2 = $ws.Range("B16").Value2
$ws.Range("C19").Value2 = $ws.Range("B19").Value2
$ws.Range("C20").Value2 = $ws.Range("B20").Value2
$ws.Range("C21").Value2 = $ws.Range("B21").Value2
$ws.Range("C22").Value2 = $ws.Range("B22").Value2
